$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the small truth table (J / K / Q(n) / Q(n-1) / ~Q(n-1)) from
#    A1:D5 down-and-right to E3:H7.
# ---------------------------------------------------------------------------
$ws.Range("A1:D5").Copy($ws.Range("E3"))

# ---------------------------------------------------------------------------
# 2) Move the lone bottom border spacer cell D6 -> H8.
# ---------------------------------------------------------------------------
$ws.Range("D6").Copy($ws.Range("H8"))

# ---------------------------------------------------------------------------
# 3) Grab the three border styles from the old JK Karnaugh-map block that are
#    still needed in the new note box (H9, E10, F10) before that block gets
#    wiped out.
# ---------------------------------------------------------------------------
$ws.Range("D7").Copy($ws.Range("H9"))
$ws.Range("A8").Copy($ws.Range("E10"))
$ws.Range("B8").Copy($ws.Range("F10"))

# ---------------------------------------------------------------------------
# 4) Clear the old JK Karnaugh-map block (now fully superseded) and the cells
#    we just vacated in step 1/2.
# ---------------------------------------------------------------------------
$ws.Range("A1:D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("A7:D15").Clear()
$ws.Range("E8:G8").Clear()

# Content we copied into E10/F10 came along with old text - strip it, the
# new note box only carries text in E9.
$ws.Range("E10:F10").ClearContents()

# ---------------------------------------------------------------------------
# 5) Extend the blank bordered note box to its full new footprint (E9:J17),
#    reusing the plain bordered-cell style already used at A16 (unaffected
#    by this edit) as the stamp source.
# ---------------------------------------------------------------------------
$ws.Range("A16").Copy($ws.Range("F9:G9"))
$ws.Range("A16").Copy($ws.Range("I9:J10"))
$ws.Range("A16").Copy($ws.Range("G10:H10"))
$ws.Range("A16").Copy($ws.Range("E11:J15"))
$ws.Range("A16").Copy($ws.Range("I16:J17"))

# ---------------------------------------------------------------------------
# 6) Style E9 like the rest of the note box and give it the new text.
# ---------------------------------------------------------------------------
$ws.Range("A16").Copy($ws.Range("E9"))
$ws.Range("E9").Value = "Zmiany stanu następują, gdy na wejściu zegarowym jest stan niski."

# ---------------------------------------------------------------------------
# 7) Selection / view bookkeeping to match the saved state.
# ---------------------------------------------------------------------------
$ws.Range("G21").Select()
